# Update the raw inputs for row 10 (VPTrial) with the corrected
# experimental values; the dependent formulas in F10/G10/H10 (Sum,
# Delta, Avg Delta) recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 10300.4335332544
$ws.Range("D10").Value = 20071.046480664001

# Leave the cursor where the author ended up after the edit.
$ws.Range("F19").Select()
